# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2210"
#   "<name>_new" -> "<name>_FV2304"
# and turn the used range A1:U66 into a real Excel Table ("Table1") with an
# AutoFilter, plus freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffixedColumns = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1..10) carried the "_old" suffix -> rename to "_FV2210".
for ($i = 0; $i -lt $suffixedColumns.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($suffixedColumns[$i] + "_FV2210")
}

# Column K (11) is "diff" and is left untouched.

# Columns L..U (12..21) carried the "_new" suffix -> rename to "_FV2304".
for ($i = 0; $i -lt $suffixedColumns.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($suffixedColumns[$i] + "_FV2304")
}

# Freeze the header row (row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn A1:U66 into a real table while preserving the header row's existing
# formatting (bold / shaded / bordered) exactly as it was: stash a copy of
# the header formatting outside the table range, let Excel apply its own
# (undesired) header styling as part of table creation, then restore the
# original formatting on top of it.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A67:U67")

$headerRange.Copy()
$scratch.PasteSpecial(-4122)
$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U66"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$scratch.ClearFormats()
$scratch.ClearContents()

Write-Output "done"
